# Apply the "Week 10" log entry to the Logboek worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 26: new "Week 10" header row (bold, like the other week-header rows) ---
$ws.Range("A26").Value = "Week 10"
$ws.Range("A26").Font.Bold = $true
$ws.Range("B26").Font.Bold = $true

# --- Row 27: the log entry itself (date + description) ---
$ws.Range("A27").Value = 43571
$ws.Range("A27").NumberFormat = $ws.Range("A25").NumberFormat
$ws.Range("C27").Value = "specifieke data ophalen uit cloud firestore en proberen tonen via markers"

# --- Row 28: new interesting-link entry with a hyperlink ---
$ws.Range("Q27").Copy()
$ws.Range("Q28").PasteSpecial(-4122) | Out-Null  # xlPasteFormats - reuse existing hyperlink style
$ws.Range("Q28").Value = "https://stackoverflow.com/questions/19282948/create-an-android-location-from-a-string-array"
$ws.Hyperlinks.Add($ws.Range("Q28"), "https://stackoverflow.com/questions/19282948/create-an-android-location-from-a-string-array") | Out-Null
$ws.Range("Q28").Style = "Hyperlink"

# --- Row 27: hours spent ---
$ws.Range("B27").Value = "3 uur 10 minuten"

# --- Row 24: add the missing hours cell next to the existing "Week 9" header (reuses existing string) ---
$ws.Range("B24").Value = "1 uur 5 minuten"
$ws.Range("B24").Font.Bold = $true

$excel.CutCopyMode = $false

# --- Update the active selection to match the final workbook state ---
$ws.Range("D28").Select()

$wb.Save()
